{"js": "// Update the worksheet date header and all 100 equation answers in the\n// table (20 rows x 5 cols) to the new values from the next day's sheet.\n//\n// REPLACEMENTS[0] is the date paragraph text; REPLACEMENTS[1..100] are\n// the equation cell texts, in row-major (top-to-bottom, left-to-right)\n// table order \u2014 matching the document's single table exactly.\nconst REPLACEMENTS = [[\"2025-09-21 Sunday\", \"2025-09-22 Monday\"], [\"60-37=23\", \"10+57=67\"], [\"82-51=31\", \"69-49=20\"], [\"76+16=92\", \"55+36=91\"], [\"63-8=55\", \"30+30=60\"], [\"23-2=21\", \"76+23=99\"], [\"41-14=27\", \"23+72=95\"], [\"43+32=75\", \"54-39=15\"], [\"18+34=52\", \"94-21=73\"], [\"83-12=71\", \"92-53=39\"], [\"15+79=94\", \"13+69=82\"], [\"20-3=17\", \"42+1=43\"], [\"47+27=74\", \"50+35=85\"], [\"7+31=38\", \"97-65=32\"], [\"1+6=7\", \"36-4=32\"], [\"54+15=69\", \"63-29=34\"], [\"64+30=94\", \"8+42=50\"], [\"68-4=64\", \"23-11=12\"], [\"70-31=39\", \"57-19=38\"], [\"63+28=91\", \"41+15=56\"], [\"26+68=94\", \"77-63=14\"], [\"55-24=31\", \"25+69=94\"], [\"99-86=13\", \"90-32=58\"], [\"74+5=79\", \"97-42=55\"], [\"7+8=15\", \"64+16=80\"], [\"41-32=9\", \"69-64=5\"], [\"95-87=8\", \"10+47=57\"], [\"68-40=28\", \"12+0=12\"], [\"5+63=68\", \"25+47=72\"], [\"78-70=8\", \"77-13=64\"], [\"59-16=43\", \"62-3=59\"], [\"96-13=83\", \"79-6=73\"], [\"4+33=37\", \"70-30=40\"], [\"4+14=18\", \"36-22=14\"], [\"65-65=0\", \"17+43=60\"], [\"3+34=37\", \"58-11=47\"], [\"83-62=21\", \"52-39=13\"], [\"13+70=83\", \"5+45=50\"], [\"52-44=8\", \"0+96=96\"], [\"7+33=40\", \"29-16=13\"], [\"30-23=7\", \"59-48=11\"], [\"0+61=61\", \"87-48=39\"], [\"62-51=11\", \"87-70=17\"], [\"23+69=92\", \"76-1=75\"], [\"85-25=60\", \"46+17=63\"], [\"47+0=47\", \"44+31=75\"], [\"7+66=73\", \"70+26=96\"], [\"36+35=71\", \"57+23=80\"], [\"89-46=43\", \"36+12=48\"], [\"93-27=66\", \"41+36=77\"], [\"17+54=71\", \"64-20=44\"], [\"91-90=1\", \"90-82=8\"], [\"91+7=98\", \"69-31=38\"], [\"57-2=55\", \"24+69=93\"], [\"34+36=70\", \"81-0=81\"], [\"21+33=54\", \"39+44=83\"], [\"8+80=88\", \"62-29=33\"], [\"54+20=74\", \"21+61=82\"], [\"5+9=14\", \"85-37=48\"], [\"24+27=51\", \"28-1=27\"], [\"93-49=44\", \"98-87=11\"], [\"31-8=23\", \"77-46=31\"], [\"88-74=14\", \"77+1=78\"], [\"17-2=15\", \"64-46=18\"], [\"12+39=51\", \"43+26=69\"], [\"75-4=71\", \"8+86=94\"], [\"40-7=33\", \"29+55=84\"], [\"94-25=69\", \"81+17=98\"], [\"56-41=15\", \"4+69=73\"], [\"11+80=91\", \"37+37=74\"], [\"6+8=14\", \"32+6=38\"], [\"72-11=61\", \"57+27=84\"], [\"31+37=68\", \"75-45=30\"], [\"99-26=73\", \"15+69=84\"], [\"15-15=0\", \"43-24=19\"], [\"5+6=11\", \"15+41=56\"], [\"39+4=43\", \"64-51=13\"], [\"72-60=12\", \"90+8=98\"], [\"50+36=86\", \"33+5=38\"], [\"33+59=92\", \"33+35=68\"], [\"3+93=96\", \"99-60=39\"], [\"77-2=75\", \"60-39=21\"], [\"35-1=34\", \"29+61=90\"], [\"16+73=89\", \"30+30=60\"], [\"72+11=83\", \"32+14=46\"], [\"39+7=46\", \"47-36=11\"], [\"34-29=5\", \"76+13=89\"], [\"55+8=63\", \"77-64=13\"], [\"76-50=26\", \"82-32=50\"], [\"70-58=12\", \"66+26=92\"], [\"23-17=6\", \"55+18=73\"], [\"81-40=41\", \"31+53=84\"], [\"91-66=25\", \"51+38=89\"], [\"72+10=82\", \"66-1=65\"], [\"3+6=9\", \"96-96=0\"], [\"18-7=11\", \"50-30=20\"], [\"15+53=68\", \"15+80=95\"], [\"23+29=52\", \"33+34=67\"], [\"21+7=28\", \"49+32=81\"], [\"78-24=54\", \"0+72=72\"], [\"45+51=96\", \"22-0=22\"]];\n\nconst body = context.document.body;\n\n// 1) Update the date paragraph (first paragraph in the document body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst [oldDate, newDate] = REPLACEMENTS[0];\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text.trim() === oldDate) {\n  // Replace via a search-and-insertText so run formatting (font/size) is kept.\n  const hits = dateParagraph.search(oldDate, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(newDate, Word.InsertLocation.replace);\n  } else {\n    dateParagraph.insertText(newDate, Word.InsertLocation.replace);\n  }\n} else {\n  // Fallback: search the whole body for the old date text.\n  const bodyHits = body.search(oldDate, { matchCase: true });\n  bodyHits.load(\"items\");\n  await context.sync();\n  if (bodyHits.items.length > 0) {\n    bodyHits.items[0].insertText(newDate, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// 2) Update every equation cell in the (single) table, in row-major order.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst tableRows = table.rows;\ntableRows.load(\"items\");\nawait context.sync();\n\nconst eqReplacements = REPLACEMENTS.slice(1);\nlet idx = 0;\nfor (const row of tableRows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (const cell of cells.items) {\n    if (idx >= eqReplacements.length) break;\n    const [oldEq, newEq] = eqReplacements[idx];\n    const cellBody = cell.body;\n    // Search for the exact old equation text inside this cell and replace\n    // just that range, so the surrounding run/paragraph formatting\n    // (font, size, paragraph alignment) is preserved.\n    const cellHits = cellBody.search(oldEq, { matchCase: true });\n    cellHits.load(\"items\");\n    await context.sync();\n    if (cellHits.items.length > 0) {\n      cellHits.items[0].insertText(newEq, Word.InsertLocation.replace);\n    } else {\n      // Fallback (should not normally trigger): replace the whole cell text.\n      cellBody.insertText(newEq, Word.InsertLocation.replace);\n    }\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date header and all 100 equation answers in the\n# table (20 rows x 5 cols) to the new values from the next day's sheet.\n# Each (Old, New) pair below is unique across the whole document body,\n# so a plain Find/ReplaceAll on $d.Content for each pair is unambiguous\n# and preserves each run's existing formatting (font/size/alignment).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2025-09-21 Sunday'; New = '2025-09-22 Monday' }\n    @{ Old = '60-37=23'; New = '10+57=67' }\n    @{ Old = '82-51=31'; New = '69-49=20' }\n    @{ Old = '76+16=92'; New = '55+36=91' }\n    @{ Old = '63-8=55'; New = '30+30=60' }\n    @{ Old = '23-2=21'; New = '76+23=99' }\n    @{ Old = '41-14=27'; New = '23+72=95' }\n    @{ Old = '43+32=75'; New = '54-39=15' }\n    @{ Old = '18+34=52'; New = '94-21=73' }\n    @{ Old = '83-12=71'; New = '92-53=39' }\n    @{ Old = '15+79=94'; New = '13+69=82' }\n    @{ Old = '20-3=17'; New = '42+1=43' }\n    @{ Old = '47+27=74'; New = '50+35=85' }\n    @{ Old = '7+31=38'; New = '97-65=32' }\n    @{ Old = '1+6=7'; New = '36-4=32' }\n    @{ Old = '54+15=69'; New = '63-29=34' }\n    @{ Old = '64+30=94'; New = '8+42=50' }\n    @{ Old = '68-4=64'; New = '23-11=12' }\n    @{ Old = '70-31=39'; New = '57-19=38' }\n    @{ Old = '63+28=91'; New = '41+15=56' }\n    @{ Old = '26+68=94'; New = '77-63=14' }\n    @{ Old = '55-24=31'; New = '25+69=94' }\n    @{ Old = '99-86=13'; New = '90-32=58' }\n    @{ Old = '74+5=79'; New = '97-42=55' }\n    @{ Old = '7+8=15'; New = '64+16=80' }\n    @{ Old = '41-32=9'; New = '69-64=5' }\n    @{ Old = '95-87=8'; New = '10+47=57' }\n    @{ Old = '68-40=28'; New = '12+0=12' }\n    @{ Old = '5+63=68'; New = '25+47=72' }\n    @{ Old = '78-70=8'; New = '77-13=64' }\n    @{ Old = '59-16=43'; New = '62-3=59' }\n    @{ Old = '96-13=83'; New = '79-6=73' }\n    @{ Old = '4+33=37'; New = '70-30=40' }\n    @{ Old = '4+14=18'; New = '36-22=14' }\n    @{ Old = '65-65=0'; New = '17+43=60' }\n    @{ Old = '3+34=37'; New = '58-11=47' }\n    @{ Old = '83-62=21'; New = '52-39=13' }\n    @{ Old = '13+70=83'; New = '5+45=50' }\n    @{ Old = '52-44=8'; New = '0+96=96' }\n    @{ Old = '7+33=40'; New = '29-16=13' }\n    @{ Old = '30-23=7'; New = '59-48=11' }\n    @{ Old = '0+61=61'; New = '87-48=39' }\n    @{ Old = '62-51=11'; New = '87-70=17' }\n    @{ Old = '23+69=92'; New = '76-1=75' }\n    @{ Old = '85-25=60'; New = '46+17=63' }\n    @{ Old = '47+0=47'; New = '44+31=75' }\n    @{ Old = '7+66=73'; New = '70+26=96' }\n    @{ Old = '36+35=71'; New = '57+23=80' }\n    @{ Old = '89-46=43'; New = '36+12=48' }\n    @{ Old = '93-27=66'; New = '41+36=77' }\n    @{ Old = '17+54=71'; New = '64-20=44' }\n    @{ Old = '91-90=1'; New = '90-82=8' }\n    @{ Old = '91+7=98'; New = '69-31=38' }\n    @{ Old = '57-2=55'; New = '24+69=93' }\n    @{ Old = '34+36=70'; New = '81-0=81' }\n    @{ Old = '21+33=54'; New = '39+44=83' }\n    @{ Old = '8+80=88'; New = '62-29=33' }\n    @{ Old = '54+20=74'; New = '21+61=82' }\n    @{ Old = '5+9=14'; New = '85-37=48' }\n    @{ Old = '24+27=51'; New = '28-1=27' }\n    @{ Old = '93-49=44'; New = '98-87=11' }\n    @{ Old = '31-8=23'; New = '77-46=31' }\n    @{ Old = '88-74=14'; New = '77+1=78' }\n    @{ Old = '17-2=15'; New = '64-46=18' }\n    @{ Old = '12+39=51'; New = '43+26=69' }\n    @{ Old = '75-4=71'; New = '8+86=94' }\n    @{ Old = '40-7=33'; New = '29+55=84' }\n    @{ Old = '94-25=69'; New = '81+17=98' }\n    @{ Old = '56-41=15'; New = '4+69=73' }\n    @{ Old = '11+80=91'; New = '37+37=74' }\n    @{ Old = '6+8=14'; New = '32+6=38' }\n    @{ Old = '72-11=61'; New = '57+27=84' }\n    @{ Old = '31+37=68'; New = '75-45=30' }\n    @{ Old = '99-26=73'; New = '15+69=84' }\n    @{ Old = '15-15=0'; New = '43-24=19' }\n    @{ Old = '5+6=11'; New = '15+41=56' }\n    @{ Old = '39+4=43'; New = '64-51=13' }\n    @{ Old = '72-60=12'; New = '90+8=98' }\n    @{ Old = '50+36=86'; New = '33+5=38' }\n    @{ Old = '33+59=92'; New = '33+35=68' }\n    @{ Old = '3+93=96'; New = '99-60=39' }\n    @{ Old = '77-2=75'; New = '60-39=21' }\n    @{ Old = '35-1=34'; New = '29+61=90' }\n    @{ Old = '16+73=89'; New = '30+30=60' }\n    @{ Old = '72+11=83'; New = '32+14=46' }\n    @{ Old = '39+7=46'; New = '47-36=11' }\n    @{ Old = '34-29=5'; New = '76+13=89' }\n    @{ Old = '55+8=63'; New = '77-64=13' }\n    @{ Old = '76-50=26'; New = '82-32=50' }\n    @{ Old = '70-58=12'; New = '66+26=92' }\n    @{ Old = '23-17=6'; New = '55+18=73' }\n    @{ Old = '81-40=41'; New = '31+53=84' }\n    @{ Old = '91-66=25'; New = '51+38=89' }\n    @{ Old = '72+10=82'; New = '66-1=65' }\n    @{ Old = '3+6=9'; New = '96-96=0' }\n    @{ Old = '18-7=11'; New = '50-30=20' }\n    @{ Old = '15+53=68'; New = '15+80=95' }\n    @{ Old = '23+29=52'; New = '33+34=67' }\n    @{ Old = '21+7=28'; New = '49+32=81' }\n    @{ Old = '78-24=54'; New = '0+72=72' }\n    @{ Old = '45+51=96'; New = '22-0=22' }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null   # 2 = wdReplaceAll\n}\n\n"}
